# "stable frontend version 1.0" -- append a small "todo" list to the bottom
# of the progress tracker (Sheet1, column C), then leave the selection on
# the last entered cell (mirrors what happens when a user types the rows in
# one after another and the cursor lands on the final cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Leave a blank gap (rows 42-43) below the existing table, matching the
# author's sheet, then add the new "todo" section in column C.
$ws.Range("C44").Value = "todo"
$ws.Range("C45").Value = "1.php socket扩容！！！"
$ws.Range("C46").Value = "2. tranfer 重写"
$ws.Range("C47").Value = "3. train inquire 阶梯表"
$ws.Range("C48").Value = "4. 购买车票 今日以前"

# Scroll/select so the view ends up parked on the newly-added last row,
# same as the saved workbook (topLeftCell A28 / selection C48).
$ws.Range("C48").Select()
try {
    $excel.ActiveWindow.ScrollRow = 28
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Window-scroll state isn't always available in headless hosts; the
    # cell selection above is the part that actually persists to the file.
}
